$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.431.48"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.601.62"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "212.36"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("D9").Value = "0.0606"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "19.25"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("D11").Value = "0.0857"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "1.826.55"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "1.601.16"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "4.00"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "0.505"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").Value = "63.66"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "26.412.09"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "232.21"
$ws.Range("E18").Value = "  +7.52%  "
$ws.Range("D19").Value = "7.69"
$ws.Range("E19").Value = "  +4.21%  "
$ws.Range("D20").Value = "0.0₃0724"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "4.26"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "2.16"
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("D24").Value = "8.94"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").Value = "146.66"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "6.99"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").Value = "15.48"
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "1.510.56"
$ws.Range("E32").Value = "  +5.93%  "
$ws.Range("D33").Value = "3.22"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").Value = "2.94"
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").Value = "0.566"
$ws.Range("E37").Value = "  -3.61%  "
$ws.Range("D38").Value = "0.0165"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("D39").Value = "0.821"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("D40").Value = "5.80"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "0.950"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("D43").Value = "2.18"
$ws.Range("E43").Value = "  +2.09%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.739.60"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "0.760"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").Value = "60.70"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "88.76"
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("D48").Value = "1.49"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").Value = "0.0960"
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0959"
$ws.Range("E51").Value = "  -8.57%  "
